$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Header row (B1:DL1) ---
$headers = @(
    "Archaeology", "Archaeology-Greek_and_Egyptian", "Archaeology-Medieval", "Archaeology-Mixed", "Archaeology-Other", "Archaeology-Prehistory", "Archaeology-Roman", "Arts-Ceramics", "Arts-Costume_and_textiles", "Arts-Crafts", "Arts-Design", "Arts-Fine_and_decorative_arts", "Arts-Glass", "Arts-Literature", "Arts-Music", "Arts-Other", "Arts-Photography", "Belief_and_identity", "Belief_and_identity-Church_treasuries", "Belief_and_identity-Ethnic_group", "Belief_and_identity-Freemasons", "Belief_and_identity-Other", "Belief_and_identity-Religion", "Belief_and_identity-Religious_buildings", "Buildings-Civic", "Buildings-Houses-Large_houses", "Buildings-Houses-Medium_houses", "Buildings-Houses-Small_houses", "Buildings-Other", "Buildings-Palace", "Buildings-Penal", "Buildings-School", "Buildings-Shops", "Communications", "Communications-Other", "Communications-Post", "Communications-Radio", "Food_and_drink", "Industry_and_manufacture-Clocks_and_watches", "Industry_and_manufacture-Industrial_life", "Industry_and_manufacture-Metals", "Industry_and_manufacture-Mining_and_quarrying", "Industry_and_manufacture-Mixed", "Industry_and_manufacture-Other", "Industry_and_manufacture-Potteries", "Industry_and_manufacture-Print", "Industry_and_manufacture-Steam_and_engines", "Industry_and_manufacture-Textiles", "Leisure_and_sport-Cricket", "Leisure_and_sport-Fairgrounds_and_amusements", "Leisure_and_sport-Film_Cinema_and_TV", "Leisure_and_sport-Other", "Leisure_and_sport-Rugby_and_football", "Leisure_and_sport-Toys_and_models", "Local_Histories", "Medicine_and_health-Hospital", "Medicine_and_health-Other", "Medicine_and_health-Professional_association", "Mixed-Bygones", "Mixed-Encyclopaedic", "Mixed-Other", "Natural_world", "Natural_world-Dinosaurs", "Natural_world-Fossils", "Natural_world-Geology", "Natural_world-Herbaria_and_gardening", "Natural_world-Mixed", "Natural_world-Other", "Natural_world-Zoology", "Other", "Personality-Art", "Personality-Explorer", "Personality-Literary", "Personality-Music", "Personality-Other", "Personality-Political", "Personality-Religious", "Personality-Scientific", "Rural_Industry-Farming", "Rural_Industry-Forges", "Rural_Industry-Other", "Rural_Industry-Rural_life", "Rural_Industry-Textiles", "Rural_Industry-Watermills", "Rural_Industry-Windmills", "Science_and_technology-Computing_and_gaming", "Science_and_technology-Other", "Sea_and_seafaring-Boats_and_ships", "Sea_and_seafaring-Fishing", "Sea_and_seafaring-Lighthouses", "Sea_and_seafaring-Mixed", "Sea_and_seafaring-Other", "Services-Fire", "Services-Other", "Services-Police", "Services-RNLI", "Transport-Aviation", "Transport-Bicycles", "Transport-Buses_and_trams", "Transport-Canals", "Transport-Cars_and_motorbikes", "Transport-Mixed", "Transport-Other", "Transport-Trains_and_railways", "Utilities-Gas_and_electricity", "Utilities-Water_and_waste", "War_and_conflict-Airforce", "War_and_conflict-Bunker", "War_and_conflict-Castles_and_forts", "War_and_conflict-Event_or_site", "War_and_conflict-Military", "War_and_conflict-Navy", "War_and_conflict-Other", "War_and_conflict-Regiment", "ROW_TOT"
)
for ($i = 0; $i -lt $headers.Length; $i++) {
    $ws.Cells.Item(1, $i + 2).Value = $headers[$i]
}

# Apply header style (bold + centered) to the newly added columns DI1:DL1
$newHeaderRange = $ws.Range($ws.Cells.Item(1, 113), $ws.Cells.Item(1, 116))
$newHeaderRange.Font.Bold = $true
$newHeaderRange.HorizontalAlignment = -4108

# --- Data rows (row 2 to row 6), columns B..DL ---
# Row 2
$row2 = @(0.024, 0.096, 0.167, 0.431, 0.048, 0.383, 0.885, 0.478, 0.407, 0.12, 0.167, 3.589, 0.215, 0.144, 0.502, 0.431, 0.191, 0, 0.144, 0.335, 0.191, 0.12, 0.478, 0.718, 0.263, 7.895, 1.388, 0.311, 0.574, 0.167, 0.311, 0.191, 0.048, 0.024, 0.12, 0.167, 0.096, 0.55, 0.144, 0.215, 0.239, 0.885, 0.191, 0.837, 0.239, 0.12, 0.215, 0.622, 0.191, 0.215, 0.478, 0.407, 0.383, 0.885, 14.713, 0.335, 0.263, 0.239, 0.311, 2.249, 1.675, 0.024, 0.096, 0.048, 0.263, 0.239, 0.478, 0.167, 0.096, 1.603, 0.359, 0.12, 1.244, 0.12, 0.622, 0.311, 0.215, 0.239, 1.196, 0.096, 0.191, 0.67, 0.167, 0.622, 0.431, 0.167, 0.311, 0.742, 0.144, 0.167, 0.789, 0.191, 0.287, 0.144, 0.502, 0.12, 0.574, 0.12, 0.311, 0.287, 1.388, 0.622, 0.239, 2.751, 0.12, 0.718, 1.579, 0.167, 1.172, 0.55, 0.909, 0.335, 0.311, 2.56, 76.31399999999999)
for ($i = 0; $i -lt $row2.Length; $i++) {
    $ws.Cells.Item(2, $i + 2).Value = $row2[$i]
}

# Row 3
$row3 = @(0, 0, 0, 0, 0, 0.024, 0, 0.024, 0, 0, 0, 0.048, 0, 0, 0.024, 0.024, 0, 0, 0, 0, 0, 0.024, 0.048, 0, 0, 0.191, 0.07199999999999999, 0.024, 0.096, 0, 0, 0.024, 0, 0, 0, 0, 0, 0, 0, 0, 0, 0, 0.024, 0, 0, 0.024, 0, 0.07199999999999999, 0, 0, 0, 0, 0.048, 0, 0.455, 0, 0, 0, 0.024, 0.048, 0.07199999999999999, 0, 0, 0, 0, 0, 0, 0, 0, 0.024, 0, 0, 0.07199999999999999, 0, 0.024, 0.12, 0.096, 0, 0.024, 0, 0, 0.048, 0, 0, 0, 0, 0, 0.096, 0, 0, 0, 0, 0, 0.024, 0.024, 0, 0.024, 0.024, 0, 0.024, 0.024, 0.024, 0, 0.096, 0.024, 0, 0, 0.024, 0, 0.12, 0.048, 0.024, 0.07199999999999999, 0.096, 2.566)
for ($i = 0; $i -lt $row3.Length; $i++) {
    $ws.Cells.Item(3, $i + 2).Value = $row3[$i]
}

# Row 4
$row4 = @(0, 0, 0.096, 0.024, 0.048, 0.215, 0.024, 0, 0.048, 0.024, 0.024, 0.742, 0, 0.048, 0.096, 0.024, 0.024, 0.024, 0.024, 0.215, 0.024, 0.024, 0.167, 0.12, 0, 1.364, 0.07199999999999999, 0.287, 0.191, 0, 0.024, 0.07199999999999999, 0.048, 0, 0.024, 0, 0.024, 0.048, 0, 0.048, 0, 0.07199999999999999, 0.048, 0.024, 0, 0, 0, 0.096, 0, 0, 0, 0.12, 0.096, 0.215, 4.593, 0.048, 0.07199999999999999, 0.048, 0.024, 0.478, 0.191, 0, 0.024, 0.048, 0.07199999999999999, 0.024, 0.07199999999999999, 0.048, 0.048, 0.191, 0.07199999999999999, 0.048, 0.502, 0, 0.239, 0.048, 0.024, 0.07199999999999999, 0.239, 0.07199999999999999, 0.048, 0.144, 0.096, 0.07199999999999999, 0, 0, 0.07199999999999999, 0.287, 0.144, 0.07199999999999999, 0.144, 0.024, 0.048, 0, 0.07199999999999999, 0, 0.048, 0, 0.048, 0.024, 0.144, 0.096, 0, 0.239, 0.024, 0, 0.096, 0.024, 0.191, 0.096, 0.096, 0.024, 0.07199999999999999, 0.263, 14.913)
for ($i = 0; $i -lt $row4.Length; $i++) {
    $ws.Cells.Item(4, $i + 2).Value = $row4[$i]
}

# Row 5
$row5 = @(0, 0.024, 0.024, 0, 0, 0, 0.048, 0.024, 0.07199999999999999, 0, 0, 0.287, 0, 0.024, 0, 0, 0, 0, 0, 0.024, 0, 0, 0, 0.096, 0, 0.431, 0.191, 0.096, 0.024, 0, 0.048, 0, 0.048, 0, 0, 0, 0.024, 0.024, 0, 0.024, 0.096, 0.239, 0.07199999999999999, 0.024, 0.024, 0, 0.024, 0.024, 0.024, 0, 0.024, 0.024, 0, 0.263, 1.459, 0.024, 0, 0, 0.024, 0.024, 0.07199999999999999, 0, 0, 0, 0.048, 0, 0, 0.024, 0, 0.07199999999999999, 0.07199999999999999, 0, 0.144, 0.024, 0.048, 0.024, 0.07199999999999999, 0, 0.07199999999999999, 0, 0.024, 0.048, 0.048, 0.024, 0.024, 0, 0, 0.12, 0.024, 0, 0.144, 0, 0.024, 0, 0.024, 0, 0.024, 0.024, 0.048, 0.024, 0.096, 0.048, 0.024, 0.287, 0, 0.024, 0.048, 0, 0.191, 0.096, 0.024, 0, 0.024, 0.167, 6.227)
for ($i = 0; $i -lt $row5.Length; $i++) {
    $ws.Cells.Item(5, $i + 2).Value = $row5[$i]
}

# Row 6
$row6 = @(0.024, 0.12, 0.287, 0.455, 0.096, 0.622, 0.957, 0.526, 0.5269999999999999, 0.144, 0.191, 4.666, 0.215, 0.216, 0.622, 0.479, 0.215, 0.024, 0.168, 0.5740000000000001, 0.215, 0.168, 0.6929999999999999, 0.9339999999999999, 0.263, 9.881, 1.723, 0.718, 0.885, 0.167, 0.383, 0.287, 0.144, 0.024, 0.144, 0.167, 0.144, 0.622, 0.144, 0.287, 0.335, 1.196, 0.335, 0.885, 0.263, 0.144, 0.239, 0.8139999999999999, 0.215, 0.215, 0.502, 0.5509999999999999, 0.527, 1.363, 21.22, 0.407, 0.335, 0.287, 0.383, 2.799, 2.01, 0.024, 0.12, 0.096, 0.383, 0.263, 0.5499999999999999, 0.239, 0.144, 1.89, 0.503, 0.168, 1.962, 0.144, 0.9329999999999999, 0.503, 0.407, 0.311, 1.531, 0.168, 0.263, 0.91, 0.311, 0.718, 0.455, 0.167, 0.383, 1.245, 0.312, 0.239, 1.077, 0.215, 0.359, 0.168, 0.622, 0.12, 0.6699999999999999, 0.168, 0.407, 0.359, 1.652, 0.79, 0.263, 3.373, 0.168, 0.742, 1.723, 0.215, 1.554, 0.8620000000000001, 1.077, 0.383, 0.479, 3.086, 100.02)
for ($i = 0; $i -lt $row6.Length; $i++) {
    $ws.Cells.Item(6, $i + 2).Value = $row6[$i]
}
